$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.169881224632263
$ws.Range("B1").Value = 2.137012720108032
$ws.Range("C1").Value = 3.09444785118103
$ws.Range("D1").Value = 3.620275259017944
$ws.Range("E1").Value = 1.4419264793396
